$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Export as TSV")
$sheet3 = $wb.Worksheets.Item("assay_type list")

# Remove the "Publication" entry (A1) from the assay_type list sheet, keeping
# only "publication" (lowercase), shifted up into A1.
$sheet3.Range("A1").Value = "publication"
$sheet3.Range("A2").ClearContents()

# Update the data validation on column B of "Export as TSV" to only reference
# the single remaining list entry.
$sheet1.Range("B2:B1048576").Validation.Delete()
$sheet1.Range("B2:B1048576").Validation.Add(3, 1, 1, "='assay_type list'!`$A`$1:`$A`$1")
$sheet1.Range("B2:B1048576").Validation.ErrorTitle = "Value must come from list"
$sheet1.Range("B2:B1048576").Validation.ErrorMessage = "Value must be one of: publication."
$sheet1.Range("B2:B1048576").Validation.ShowInput = $true
$sheet1.Range("B2:B1048576").Validation.ShowError = $true
